$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Title row (row 1): alignment changes from centered to left-aligned.
# ---------------------------------------------------------------------------
$ws.Range("A1:F1").HorizontalAlignment = -4131   # xlLeft

# ---------------------------------------------------------------------------
# 2. New sub-header row 2: "After 1000 iterations" merged across B2:C2,
#    centered, using the default (non-bold, non-title) font.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "After 1000 iterations"
$ws.Range("B2:C2").Merge()
$ws.Range("B2:C2").HorizontalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 3. Fill in the "After 1000 iterations" Training/Testing error results for
#    each sigma value (rows 4-13, columns B and C).
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = 0.00624294896705
$ws.Range("C4").Value = 0.0963513527258
$ws.Range("B5").Value = 0.0418015127285
$ws.Range("C5").Value = 0.0746117360473
$ws.Range("B6").Value = 0.053880318463
$ws.Range("C6").Value = 0.0667372576432
$ws.Range("B7").Value = 0.0595791092407
$ws.Range("C7").Value = 0.0630058579812
$ws.Range("B8").Value = 0.0613869979149
$ws.Range("C8").Value = 0.064432628574
$ws.Range("B9").Value = 0.0677665134962
$ws.Range("C9").Value = 0.0671847768536
$ws.Range("B10").Value = 0.0751723451117
$ws.Range("C10").Value = 0.0702082150911
$ws.Range("B11").Value = 0.0787587990829
$ws.Range("C11").Value = 0.0704839293719
$ws.Range("B12").Value = 0.0810244640352
$ws.Range("C12").Value = 0.0696689145693
$ws.Range("B13").Value = 0.0819406404426
$ws.Range("C13").Value = 0.0691972653316

# ---------------------------------------------------------------------------
# 4. Column widths widen slightly now that columns B and C hold data.
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 12.5
$ws.Columns("B").ColumnWidth = 13.333333333333334
$ws.Columns("C").ColumnWidth = 12.666666666666666

# ---------------------------------------------------------------------------
# 5. Selection moves to the last filled cell, C13.
# ---------------------------------------------------------------------------
$ws.Range("C13").Select()
